$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for even_MAG-GUT38177.fa (row 13), shifting subsequent rows up.
$ws.Rows.Item(13).Delete()
